$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
  @(8,9),
  @(2,3),
  @(1,6),
  @(1,5),
  @(1,7),
  @(1,8),
  @(1,5),
  @(1,6),
  @(1,6),
  @(1,6),
  @(1,6),
  @(1,6),
  @(1,5),
  @(8,9),
  @(1,4),
  @(6,7),
  @(1,2),
  @(8,8),
  @(1,2),
  @(6,8),
  @(9,9),
  @(8,8),
  @(7,8),
  @(9,9),
  @(9,9),
  @(6,6),
  @(8,8),
  @(7,8),
  @(8,8),
  @(5,5),
  @(1,1),
  @(9,9),
  @(7,7),
  @(7,8),
  @(8,8),
  @(6,6),
  @(7,7)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 9).Value = $data[$i][0]
  $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
